$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.312.49"
$ws.Range("E2").Value = "  -2.51%  "
$ws.Range("D3").Value = "2.573.47"
$ws.Range("E3").Value = "  -2.92%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.18"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.22"
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.75"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0996"
$ws.Range("E10").Value = "  -3.57%  "
$ws.Range("E11").Value = "  +2.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.330"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").Value = "3.026.79"
$ws.Range("E13").Value = "  -3.09%  "
$ws.Range("D14").Value = "58.213.06"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.47"
$ws.Range("E15").Value = "  -3.84%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.569.14"
$ws.Range("E16").Value = "  -3.91%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.44"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "333.18"
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.98"
$ws.Range("E20").Value = "  -2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.11"
$ws.Range("E21").Value = "  -3.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.39"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.417"
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.157"
$ws.Range("E26").Value = "  -5.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.99"
$ws.Range("E27").Value = "  -4.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D29").Value = "0.0₃0727"
$ws.Range("E29").Value = "  -3.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.65"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.11"
$ws.Range("E31").Value = "  +3.27%  "
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "36.93"
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.842"
$ws.Range("E36").Value = "  +2.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.08"
$ws.Range("E37").Value = "  -5.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.812"
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("E39").Value = "  -4.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.56"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "278.33"
$ws.Range("E41").Value = "  -6.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.62"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.586"
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0938"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.35"
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0225"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").Value = "1.903.51"
$ws.Range("E49").Value = "  -3.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.74"
$ws.Range("E50").Value = "  -3.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.35"
$ws.Range("E51").Value = "  -4.99%  "
